$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a serial date value that was bumped by one day
# (45188 -> 45189) for every data row (rows 2 through 205).
$ws.Range("C2:C205").Value = 45189
